$d = $word.ActiveDocument
$t = $d.Tables.Item(7)
$cell = $t.Cell(2, 2)
$p1 = $cell.Range.Paragraphs.Item(1)
$nl = [char]13
$texts = @(
  "Began implementing shape context",
  "Sampled points along edge of image and calculate distance between point pairs. Normalize distance",
  "Started calculating average vector between point pairs to approximate tangent line at each point",
  "Angle between point pairs to be found using tangent line rather than x-axis to be invariant to rotation"
)
$full = [string]::Join($nl, $texts)
$p1.Range.Text = $full
